$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.424.07'
$ws.Range('E2').Value = '  -0.09%  '
$ws.Range('D3').Value = '1.904.59'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('E5').Value = '  +10.64%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '246.74'
$ws.Range('E6').Value = '  +0.63%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '40.69'
$ws.Range('E8').Value = '  -3.04%  '
$ws.Range('E9').Value = '  +3.50%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '52.55'
$ws.Range('E10').Value = '  +7.91%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0726'
$ws.Range('E11').Value = '  +3.21%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0990'
$ws.Range('E12').Value = '  -0.60%  '
$ws.Range('D13').Value = '2.181.60'
$ws.Range('E13').Value = '  +0.06%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '12.62'
$ws.Range('E14').Value = '  +2.34%  '
$ws.Range('E15').Value = '  +2.95%  '
$ws.Range('D16').Value = '1.911.84'
$ws.Range('E16').Value = '  +0.56%  '
$ws.Range('E17').Value = '  +0.72%  '
$ws.Range('D18').Value = '35.398.19'
$ws.Range('E18').Value = '  -0.28%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.03'
$ws.Range('E19').Value = '  +1.74%  '
$ws.Range('D20').Value = '0.0₃0824'
$ws.Range('E20').Value = '  +0.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '241.78'
$ws.Range('E21').Value = '  -0.54%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.93'
$ws.Range('E22').Value = '  +2.78%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.10'
$ws.Range('E23').Value = '  +5.05%  '
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('E25').Value = '  +0.86%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.32'
$ws.Range('E26').Value = '  +7.22%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '169.15'
$ws.Range('E27').Value = '  -1.87%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.66'
$ws.Range('E28').Value = '  +1.51%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.89'
$ws.Range('E29').Value = '  +5.26%  '
$ws.Range('E30').Value = '  +4.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.23'
$ws.Range('E32').Value = '  +3.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0575'
$ws.Range('E33').Value = '  +0.54%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.88'
$ws.Range('E34').Value = '  +7.41%  '
$ws.Range('B35').Value = 'BinanceUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.01'
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('B36').Value = 'InternetComputer(DFINITY)'
$ws.Range('C36').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.18'
$ws.Range('E36').Value = '  +0.61%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.918'
$ws.Range('E37').Value = '  -5.86%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.48'
$ws.Range('E38').Value = '  +8.93%  '
$ws.Range('E39').Value = '  +0.07%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '96.93'
$ws.Range('E40').Value = '  +6.80%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.11'
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '16.59'
$ws.Range('E42').Value = '  +5.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0657'
$ws.Range('E43').Value = '  +5.00%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0209'
$ws.Range('E44').Value = '  +1.75%  '
$ws.Range('D45').Value = '1.357.29'
$ws.Range('E45').Value = '  +0.53%  '
$ws.Range('E46').Value = '  +2.16%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '46.32'
$ws.Range('E47').Value = '  -8.39%  '
$ws.Range('E48').Value = '  +0.38%  '
$ws.Range('E49').Value = '  +1.10%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '12.31'
$ws.Range('E50').Value = '  -4.97%  '
$ws.Range('E51').Value = '  -1.31%  '
